# CA test updated and xml added
#
# CA_FIL sheet ("Scenario2"/"Scenario3" rows removed, "Scenario4" row kept as the
# sole remaining scenario, with its Location changed to "Null"), and the author's
# last-saved UI state (active tab / selected cell) moved onto CA_FIL.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("CA_FIL")

# Drop the three scenario rows that are no longer needed (old rows 3,4,5 -
# Scenario2/Scenario3/Scenario4); this leaves the header row and the old
# Scenario1 row, which we immediately overwrite below with the data that the
# surviving "Scenario4" row should hold.
$ws3.Rows("3:5").Delete()

$ws3.Range("A2").Value = "Scenario4"
$ws3.Range("B2").Value = "CA"
$ws3.Range("C2").Value = "Ontario"
$ws3.Range("D2").Value = "North York"
$ws3.Range("E2").Value = "Null"

# Match the saved workbook UI state: CA_FIL tab active, E2 selected there.
$ws3.Activate()
$ws3.Range("E2").Select()
